$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 12: new finished time-log entry ("Code cleanup") ---
$ws.Range("A12").Value = 44502
$ws.Range("B12").Value = 0.3888888888888889
$ws.Range("C12").Value = 0.42708333333333331
$ws.Range("D12").NumberFormat = $ws.Range("D11").NumberFormat
$ws.Range("D12").Formula = "=C12-B12"
$ws.Range("E12").Value = "Code cleanup"
$ws.Range("F12").Value = "Cleaning up code for presentation."

# --- Row 13: started next entry ("Presentation Prep"), no end time yet ---
$ws.Range("B13").Value = 0.79166666666666663
$ws.Range("E13").Value = "Presentation Prep"
$ws.Range("F13").Value = "Preparing for Phase 2 Presentation"

# --- Update selection to match the author's last-selected cell ---
$ws.Range("E14").Select()
